# Populate Sheet2 with the per-year sbmsy ratio data (2003-2012) and a
# trailing average, then add an XY-scatter chart with a linear trendline
# plotting year (D) vs ratio (E).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: raw numerator / denominator, ratio formula, year, ratio again ---
$numerators = @(110101, 110655, 109212, 105012, 98562, 106245, 110282, 112767, 113844, 111833)
$years      = @(2012, 2011, 2010, 2009, 2008, 2007, 2006, 2005, 2004, 2003)

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 1
    $ws2.Cells.Item($r, 1).Value = $numerators[$i]
    $ws2.Cells.Item($r, 2).Value = 49680
    $ws2.Range("C$r").Formula = "=A$r/B$r"
    $ws2.Cells.Item($r, 4).Value = $years[$i]
}

$excel.Calculate()

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 1
    $ws2.Cells.Item($r, 5).Value = $ws2.Cells.Item($r, 3).Value2
}

$ws2.Range("E11").Formula = "=AVERAGE(E1:E10)"
$ws2.Range("A1:A10").NumberFormat = "#,##0"

$ws2.Range("E12").Select()

# --- Chart: XY scatter of year (x) vs ratio (y) with a linear trendline ---
$chartObj = $ws2.ChartObjects().Add(228600, 68580, 4305300, 2971800)
$chart = $chartObj.Chart
$chart.ChartType = -4169
$series = $chart.SeriesCollection().NewSeries()
$series.XValues = $ws2.Range("D1:D10")
$series.Values = $ws2.Range("E1:E10")
$trendline = $series.Trendlines().Add(-4132)
$trendline.DisplayRSquared = $true
$trendline.DisplayEquation = $false

# --- Restore Sheet1 as the active tab with the new selection; Sheet2 keeps
#     its own cursor at E12 for when it is next activated. ---
$ws1.Select()
$ws1.Range("E5:G9").Select()

$wb.Save()
